$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("policies")

$ws.Range("D4").Value = "Grid_MS_congestion_price_eurpkWh"
$ws.Range("D5").Value = "Grid_MS_congestion_threshold_fr"
$ws.Range("D6").Value = "Grid_MS_congestion_pricing_consumption_eurpkWh"
$ws.Range("D7").Value = "Grid_MS_congestion_pricing_production_eurpkWh"
$ws.Range("E7").Value = "TRUE"
$ws.Range("D8").Value = "Fixed_electricity_price_eurpkWh"
$ws.Range("D9").Value = "Fixed_heat_price_eurpkWh"
$ws.Range("D10").Value = "Fixed_methane_price_eurpkWh"
$ws.Range("D11").Value = "Fixed_hydrogen_price_eurpkWh"
$ws.Range("D12").Value = "Energy_supplier_electricity_price_margin_eurpkWh"
$ws.Range("D13").Value = "Fixed_electricity_feed_in_tariff_eurpkWh"
$ws.Range("E13").Value = "0.1"
$ws.Range("D14").Value = "Fixed_diesel_price_eurpkWh"
$ws.Range("E14").Value = "0.2"
$ws.Range("D15").Value = "Time_buffer_for_spread_charging_min"
